$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct cell M20: was "NAY" (erroneous), should be "NAT" matching row 21
$ws.Range("M20").Value = "NAT"

# Widen column F (genotype1) to fit content
$ws.Columns.Item(6).ColumnWidth = 13.17

# Update the active cell selection
$ws.Range("M7").Select()
